# "same elements array researches has been added"
#
# Adds two new worksheets at the end of the workbook:
#   - same_elements  (after partly_sorted)
#   - partly_same    (after same_elements)
# and fills them with the same A1:G5 benchmark-style layout used by the
# existing sheets (straight/reversed/sorted/partly_sorted): a header row of
# array sizes in B1:G1, and one data row per element-type ("byte","int",
# "string","date") carrying elapsed-time style numbers in B:G.

$wb = $excel.ActiveWorkbook

# --- add the two new sheets at the end, in order ---------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

$wsSame = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsSame.Name = "same_elements"

$wsPartly = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsSame)
$wsPartly.Name = "partly_same"

# --- same_elements -----------------------------------------------------
$wsSame.Range("B1").Value = 5
$wsSame.Range("C1").Value = 50
$wsSame.Range("D1").Value = 500
$wsSame.Range("E1").Value = 5000
$wsSame.Range("F1").Value = 50000
$wsSame.Range("G1").Value = 500000

$wsSame.Range("A2").Value = "byte"
$wsSame.Range("B2").Value = 0
$wsSame.Range("C2").Value = 0
$wsSame.Range("D2").Value = 0.038005
$wsSame.Range("E2").Value = 3.675999
$wsSame.Range("F2").Value = 0
$wsSame.Range("G2").Value = 0

$wsSame.Range("A3").Value = "int"
$wsSame.Range("B3").Value = 0.001
$wsSame.Range("C3").Value = 0
$wsSame.Range("D3").Value = 0.039003
$wsSame.Range("E3").Value = 2.741864
$wsSame.Range("F3").Value = 0
$wsSame.Range("G3").Value = 0

$wsSame.Range("A4").Value = "string"
$wsSame.Range("B4").Value = 0
$wsSame.Range("C4").Value = 0
$wsSame.Range("D4").Value = 0.025002
$wsSame.Range("E4").Value = 1.994583
$wsSame.Range("F4").Value = 0
$wsSame.Range("G4").Value = 0

$wsSame.Range("A5").Value = "date"
$wsSame.Range("B5").Value = 0
$wsSame.Range("C5").Value = 0
$wsSame.Range("D5").Value = 0.028002
$wsSame.Range("E5").Value = 2.128223
$wsSame.Range("F5").Value = 0
$wsSame.Range("G5").Value = 0

# --- partly_same ---------------------------------------------------------
$wsPartly.Range("B1").Value = 5
$wsPartly.Range("C1").Value = 50
$wsPartly.Range("D1").Value = 500
$wsPartly.Range("E1").Value = 5000
$wsPartly.Range("F1").Value = 50000
$wsPartly.Range("G1").Value = 500000

$wsPartly.Range("A2").Value = "byte"
$wsPartly.Range("B2").Value = 0
$wsPartly.Range("C2").Value = 0
$wsPartly.Range("D2").Value = 0.004504
$wsPartly.Range("E2").Value = 0.46248
$wsPartly.Range("F2").Value = 0
$wsPartly.Range("G2").Value = 0

$wsPartly.Range("A3").Value = "int"
$wsPartly.Range("B3").Value = 0
$wsPartly.Range("C3").Value = 0
$wsPartly.Range("D3").Value = 0.016549
$wsPartly.Range("E3").Value = 0.805808
$wsPartly.Range("F3").Value = 0
$wsPartly.Range("G3").Value = 0

$wsPartly.Range("A4").Value = "string"
$wsPartly.Range("B4").Value = 0
$wsPartly.Range("C4").Value = 0
$wsPartly.Range("D4").Value = 0.01802
$wsPartly.Range("E4").Value = 2.180778
$wsPartly.Range("F4").Value = 0
$wsPartly.Range("G4").Value = 0

$wsPartly.Range("A5").Value = "date"
$wsPartly.Range("B5").Value = 0
$wsPartly.Range("C5").Value = 0.000501
$wsPartly.Range("D5").Value = 0.016521
$wsPartly.Range("E5").Value = 2.258496
$wsPartly.Range("F5").Value = 0
$wsPartly.Range("G5").Value = 0

# --- final UI state: same_elements tab active, H21 selected there ---------
$wsSame.Activate()
$wsSame.Range("H21").Select() | Out-Null
